$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 612, shifting existing rows 612.. down by one.
$ws.Rows.Item(612).Insert(-4121)

# Populate the newly inserted row 612 with the new weekly data point.
$ws.Cells.Item(612, 1).Value = 8
$ws.Cells.Item(612, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(612, 3).Value = "Coquimbo"
$ws.Cells.Item(612, 4).Value = 45267
$ws.Cells.Item(612, 5).Value = 4
$ws.Cells.Item(612, 6).Value = 100114013
$ws.Cells.Item(612, 7).Value = "Zanahoria"
$ws.Cells.Item(612, 8).Value = "Sin especificar"
$ws.Cells.Item(612, 9).Value = "Primera"
$ws.Cells.Item(612, 10).Value = 460
$ws.Cells.Item(612, 11).Value = 5500
$ws.Cells.Item(612, 12).Value = 6000
$ws.Cells.Item(612, 13).Value = 5750
$ws.Cells.Item(612, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(612, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(612, 16).Value = 288
$ws.Cells.Item(612, 17).Value = 20
$ws.Cells.Item(612, 18).Value = "Hortaliza"

# Copy the date-number-format style from the row above onto the new D612 cell.
$ws.Cells.Item(611, 4).Copy()
$ws.Cells.Item(612, 4).PasteSpecial(-4122)
